$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 2
$ws.Range("G12").Formula = '=B12*$K$6'

$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)

$ws.Range("A13").Value = 43504
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Implementation"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "Create Layout, Main Controller"
$ws.Range("F13").Value = "Need to ask about that auto complete html header thing?"

$ws.Range("A14").Value = 43504
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Implementation"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "Create Layout, Main Controller"
$ws.Range("F14").Value = "Created structure ^ stitch in time. Content still basic."

$ws.Columns.Item(6).ColumnWidth = 57.3

$ws.Range("F14").Select() | Out-Null
